$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (JOBY)
$ws.Range("D2").Value = 13.83
$ws.Range("E2").Value = 29.3
$ws.Range("F2").Value = -0.43
$ws.Range("G2").Value = 20
$ws.Range("K2").Value = 55
$ws.Range("N2").Value = 66.04328690552585

# Row 3 (ACHR)
$ws.Range("D3").Value = 7.65
$ws.Range("E3").Value = 33.9
$ws.Range("F3").Value = 2.82
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 51.2
$ws.Range("N3").Value = 66.04328690552585
